$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Date Colours")

# -----------------------------------------------------------------
# Sheet1 / Table1 : add three new outbreak-path rows (30-32)
# -----------------------------------------------------------------

# Row 30 - 2021-07-07 : Q20 w -> Q25
$ws1.Range("A30").Value = 44384
$ws1.Range("A30").NumberFormat = "d-mmm"
$ws1.Range("B30").Value = "Q20 w"
$ws1.Range("C30").Value = "Q25"
$ws1.Range("D30").Value = "Queensland"
$ws1.Range("G30").Value = "Alpha (B.1.1.7)"
$ws1.Range("H30").Value = "Isolated"

# Row 31 - 2021-07-08 : Q14 w -> Q26
$ws1.Range("A31").Value = 44385
$ws1.Range("A31").NumberFormat = "d-mmm"
$ws1.Range("B31").Value = "Q14 w"
$ws1.Range("C31").Value = "Q26"
$ws1.Range("D31").Value = "Queensland"
$ws1.Range("G31").Value = "Alpha (B.1.1.7)"
$ws1.Range("H31").Value = "Isolated"

# Row 32 - 2021-07-08 : Q14 w -> Q27
$ws1.Range("A32").Value = 44385
$ws1.Range("A32").NumberFormat = "d-mmm"
$ws1.Range("B32").Value = "Q14 w"
$ws1.Range("C32").Value = "Q27"
$ws1.Range("D32").Value = "Queensland"
$ws1.Range("G32").Value = "Alpha (B.1.1.7)"
$ws1.Range("H32").Value = "Isolated"

# Grow the structured table so the new rows participate in it (and the
# autoFilter range grows with it) - mirrors Table1 ref going A1:H29 -> A1:H32
$table1 = $ws1.ListObjects.Item("Table1")
$table1.Resize($ws1.Range("A1:H32"))

# Reflect the new selection left behind on Sheet1 after the edit
$ws1.Range("H31:H32").Select()

# -----------------------------------------------------------------
# Sheet2 ("Date Colours") : drop the stray F2:V2 helper cells and
# refresh the Colour Code dbRaevn gradient (column B, rows 2-20)
# -----------------------------------------------------------------

$ws2.Range("F2:V2").ClearContents()

$ws2.Range("B2").Value  = "#faf6fa"
$ws2.Range("B3").Value  = "#f5edf6"
$ws2.Range("B4").Value  = "#f0e4f1"
$ws2.Range("B5").Value  = "#ebdbec"
$ws2.Range("B6").Value  = "#e6d2e7"
$ws2.Range("B7").Value  = "#e1c9e3"
$ws2.Range("B8").Value  = "#dcc0de"
$ws2.Range("B9").Value  = "#d7b7d9"
$ws2.Range("B10").Value = "#d2aed5"
$ws2.Range("B11").Value = "#cda5d0"
$ws2.Range("B12").Value = "#c89dcb"
$ws2.Range("B13").Value = "#c394c7"
$ws2.Range("B14").Value = "#be8bc2"
$ws2.Range("B15").Value = "#b983bd"
$ws2.Range("B16").Value = "#b37ab9"
$ws2.Range("B17").Value = "#ae71b4"
$ws2.Range("B18").Value = "#a969af"
$ws2.Range("B19").Value = "#a360ab"
$ws2.Range("B20").Value = "#9e57a6"
